# Update "paises.xlsx" COVID data sheet:
#  - refresh the "datos actualizados" timestamp
#  - refresh numeric case counts for a handful of countries
#  - re-rank four country pairs whose totals crossed one another,
#    which also makes the newly-overtaking country show its refreshed
#    numbers while the overtaken country keeps its (unchanged) numbers
#    one row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$Country,
        [double]$Total,
        [double]$NewCases,
        [double]$Active,
        [double]$Recovered,
        [double]$Critical,
        [double]$DeathsToday,
        [double]$Deaths
    )

    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $NewCases
    $ws.Cells.Item($Row, 4).Value = $Active
    $ws.Cells.Item($Row, 5).Value = $Recovered
    $ws.Cells.Item($Row, 6).Value = $Critical
    $ws.Cells.Item($Row, 7).Value = $DeathsToday
    $ws.Cells.Item($Row, 8).Value = $Deaths
}

# --- Timestamp (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 17:22"

# --- Straight numeric refreshes (no re-ranking) -------------------------
Set-Row 4   "Estados Unidos" 3236047 16048 1426613 1673456 0 156 135978
Set-Row 6   "India"          804861  10019 503746  279339  0 153 21776
Set-Row 11  "Reino Unido"    288133  512   0       0       0 48  44650
Set-Row 19  "Alemania"       199257  59    183600  6531    0 1   9126
Set-Row 63  "Serbia"         17728   386   13651   3707    0 18  370
Set-Row 110 "Sri Lanka"      2450    296   1980    459     0 0   11

# --- Re-ranked pairs: the overtaking country gets the refreshed figures
#     and takes the higher row; the overtaken country keeps its previous
#     (unchanged) figures and drops one row. ---------------------------

# Irak overtakes Indonesia
Set-Row 29 "Irak"      72460 2848 41380 28120 0 78 2960
Set-Row 30 "Indonesia" 72347 1611 33529 35349 0 52 3469

# Portugal overtakes Singapur
Set-Row 41 "Portugal" 45679 402 30350 13683 0 2 1646
Set-Row 42 "Singapur" 45614 191 41645 3943  0 0 26

# Moldavia overtakes Austria
Set-Row 60 "Moldavia" 18924 258 12188 6101 0 11 635
Set-Row 61 "Austria"  18709 94  16808 1195 0 0  706

# Surinam overtakes San Marino
Set-Row 152 "Surinam"    699 5 456 226 0 0 17
Set-Row 153 "San Marino" 699 0 656 1   0 0 42
